$wb = $excel.ActiveWorkbook

# Update trafo model parameters (I_E_pu and P_Fe_pu) on the "trafos" sheet
$trafos = $wb.Worksheets.Item("trafos")
$trafos.Range("G2").Value = 0.05
$trafos.Range("H2").Value = 0.001

# Move the cell selection on the trafos sheet to J22
$trafos.Activate()
$trafos.Range("J22").Select()

# Make "lines" the active sheet/tab
$lines = $wb.Worksheets.Item("lines")
$lines.Activate()
